$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '39.922.79'
$ws.Cells.Item(2,5).Value = '  +1.39%  '

$ws.Cells.Item(3,4).Value = '2.219.79'
$ws.Cells.Item(3,5).Value = '  +1.02%  '

$ws.Cells.Item(4,5).Value = '  +0.01%  '

$c = $ws.Cells.Item(5,4)
$c.NumberFormat = '@'
$c.Value = '293.26'
$ws.Cells.Item(5,5).Value = '  -0.65%  '

$c = $ws.Cells.Item(6,4)
$c.NumberFormat = '@'
$c.Value = '86.94'
$ws.Cells.Item(6,5).Value = '  +6.42%  '

$ws.Cells.Item(7,5).Value = '  +0.83%  '

$ws.Cells.Item(8,5).Value = '  -0.11%  '

$c = $ws.Cells.Item(9,4)
$c.NumberFormat = '@'
$c.Value = '0.473'
$ws.Cells.Item(9,5).Value = '  +1.41%  '

$c = $ws.Cells.Item(10,4)
$c.NumberFormat = '@'
$c.Value = '30.55'
$ws.Cells.Item(10,5).Value = '  +4.85%  '

$c = $ws.Cells.Item(11,4)
$c.NumberFormat = '@'
$c.Value = '0.0786'
$ws.Cells.Item(11,5).Value = '  +2.17%  '

$c = $ws.Cells.Item(12,4)
$c.NumberFormat = '@'
$c.Value = '47.56'
$ws.Cells.Item(12,5).Value = '  +0.85%  '

$c = $ws.Cells.Item(13,4)
$c.NumberFormat = '@'
$c.Value = '0.109'
$ws.Cells.Item(13,5).Value = '  +1.67%  '

$ws.Cells.Item(14,5).Value = '  +1.94%  '

$ws.Cells.Item(15,4).Value = '2.562.21'
$ws.Cells.Item(15,5).Value = '  +1.01%  '

$c = $ws.Cells.Item(16,4)
$c.NumberFormat = '@'
$c.Value = '14.04'
$ws.Cells.Item(16,5).Value = '  +0.57%  '

$ws.Cells.Item(17,4).Value = '2.215.14'
$ws.Cells.Item(17,5).Value = '  +0.61%  '

$c = $ws.Cells.Item(18,4)
$c.NumberFormat = '@'
$c.Value = '0.731'
$ws.Cells.Item(18,5).Value = '  +2.86%  '

$ws.Cells.Item(19,4).Value = '39.847.35'
$ws.Cells.Item(19,5).Value = '  +1.47%  '

$ws.Cells.Item(20,5).Value = '  +1.44%  '

$c = $ws.Cells.Item(21,4)
$c.NumberFormat = '@'
$c.Value = '11.30'
$ws.Cells.Item(21,5).Value = '  +9.88%  '

$c = $ws.Cells.Item(22,4)
$c.NumberFormat = '@'
$c.Value = '5.81'
$ws.Cells.Item(22,5).Value = '  +1.72%  '

$c = $ws.Cells.Item(23,4)
$c.NumberFormat = '@'
$c.Value = '65.69'
$ws.Cells.Item(23,5).Value = '  +1.46%  '

$c = $ws.Cells.Item(24,4)
$c.NumberFormat = '@'
$c.Value = '236.23'
$ws.Cells.Item(24,5).Value = '  +4.89%  '

$c = $ws.Cells.Item(25,4)
$c.NumberFormat = '@'
$c.Value = '0.999'
$ws.Cells.Item(25,5).Value = '  -0.32%  '

$ws.Cells.Item(26,5).Value = '  +2.99%  '

$ws.Cells.Item(27,5).Value = '  +2.46%  '

$c = $ws.Cells.Item(28,4)
$c.NumberFormat = '@'
$c.Value = '22.80'
$ws.Cells.Item(28,5).Value = '  +1.06%  '

$c = $ws.Cells.Item(29,4)
$c.NumberFormat = '@'
$c.Value = '2.19'
$ws.Cells.Item(29,5).Value = '  +1.09%  '

$c = $ws.Cells.Item(30,4)
$c.NumberFormat = '@'
$c.Value = '9.28'
$ws.Cells.Item(30,5).Value = '  +2.23%  '

$c = $ws.Cells.Item(31,4)
$c.NumberFormat = '@'
$c.Value = '32.88'
$ws.Cells.Item(31,5).Value = '  +3.88%  '

$c = $ws.Cells.Item(32,4)
$c.NumberFormat = '@'
$c.Value = '151.97'
$ws.Cells.Item(32,5).Value = '  +2.30%  '

$c = $ws.Cells.Item(33,4)
$c.NumberFormat = '@'
$c.Value = '0.999'
$ws.Cells.Item(33,5).Value = '  -0.18%  '

$ws.Cells.Item(34,5).Value = '  +3.15%  '

$c = $ws.Cells.Item(35,4)
$c.NumberFormat = '@'
$c.Value = '0.0720'
$ws.Cells.Item(35,5).Value = '  +3.82%  '

$ws.Cells.Item(36,5).Value = '  +1.82%  '

$c = $ws.Cells.Item(37,4)
$c.NumberFormat = '@'
$c.Value = '2.80'
$ws.Cells.Item(37,5).Value = '  +6.22%  '

$ws.Cells.Item(38,5).Value = '  +2.06%  '

$c = $ws.Cells.Item(39,4)
$c.NumberFormat = '@'
$c.Value = '15.92'
$ws.Cells.Item(39,5).Value = '  +3.40%  '

$c = $ws.Cells.Item(40,4)
$c.NumberFormat = '@'
$c.Value = '0.0992'
$ws.Cells.Item(40,5).Value = '  +3.39%  '

$c = $ws.Cells.Item(41,4)
$c.NumberFormat = '@'
$c.Value = '1.71'
$ws.Cells.Item(41,5).Value = '  +3.91%  '

$c = $ws.Cells.Item(42,4)
$c.NumberFormat = '@'
$c.Value = '3.80'
$ws.Cells.Item(42,5).Value = '  +5.43%  '

$ws.Cells.Item(43,4).Value = '2.064.19'
$ws.Cells.Item(43,5).Value = '  +8.84%  '

$ws.Cells.Item(44,2).Value = 'EnergySwap'
$ws.Cells.Item(44,3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Cells.Item(44,4)
$c.NumberFormat = '@'
$c.Value = '18.00'
$ws.Cells.Item(44,5).Value = '  +12.20%  '

$ws.Cells.Item(45,2).Value = 'VeChain'
$ws.Cells.Item(45,3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Cells.Item(45,4)
$c.NumberFormat = '@'
$c.Value = '0.0268'
$ws.Cells.Item(45,5).Value = '  +3.63%  '

$ws.Cells.Item(46,2).Value = 'FraxShare'
$ws.Cells.Item(46,3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Cells.Item(46,4)
$c.NumberFormat = '@'
$c.Value = '9.98'
$ws.Cells.Item(46,5).Value = '  +11.76%  '

$c = $ws.Cells.Item(47,4)
$c.NumberFormat = '@'
$c.Value = '2.10'
$ws.Cells.Item(47,5).Value = '  +0.92%  '

$ws.Cells.Item(48,5).Value = '  +1.22%  '

$ws.Cells.Item(49,4).Value = '2.434.71'
$ws.Cells.Item(49,5).Value = '  +1.28%  '

$c = $ws.Cells.Item(50,4)
$c.NumberFormat = '@'
$c.Value = '71.31'
$ws.Cells.Item(50,5).Value = '  -0.06%  '

$c = $ws.Cells.Item(51,4)
$c.NumberFormat = '@'
$c.Value = '89.34'
$ws.Cells.Item(51,5).Value = '  +2.60%  '
